$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their exact text representation
# (values like "316.42" look numeric and would otherwise be reinterpreted
# as floating point numbers by Excel, losing trailing zeros / exact format).
$dCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.604.19'
$ws.Range('E2').Value = '  +0.99%  '

$ws.Range('D3').Value = '1.872.02'
$ws.Range('E3').Value = '  -0.17%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '316.42'
$ws.Range('E5').Value = '  +0.54%  '

$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('D7').Value = '0.5067'
$ws.Range('E7').Value = '  -1.28%  '

$ws.Range('D8').Value = '0.3891'
$ws.Range('E8').Value = '  -1.22%  '

$ws.Range('D9').Value = '0.08368'
$ws.Range('E9').Value = '  +1.03%  '

$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '41.85'
$ws.Range('E10').Value = '  +0.07%  '

$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value = '1.103'
$ws.Range('E11').Value = '  -1.44%  '

$ws.Range('D12').Value = '6.208'
$ws.Range('E12').Value = '  -1.23%  '

$ws.Range('D13').Value = '1.876.36'
$ws.Range('E13').Value = '  +0.37%  '

$ws.Range('D14').Value = '20.36'
$ws.Range('E14').Value = '  +0.11%  '

$ws.Range('D15').Value = '7.223'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').Value = '1.014'
$ws.Range('E16').Value = '  +0.16%  '

$ws.Range('D17').Value = '0.00001103'
$ws.Range('E17').Value = '  -0.18%  '

$ws.Range('D18').Value = '90.96'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('D19').Value = '0.06706'
$ws.Range('E19').Value = '  -0.13%  '

$ws.Range('D20').Value = '17.70'
$ws.Range('E20').Value = '  -0.19%  '

$ws.Range('D21').Value = '1.011'
$ws.Range('E21').Value = '  +0.10%  '

$ws.Range('D22').Value = '5.923'
$ws.Range('E22').Value = '  -1.27%  '

$ws.Range('D23').Value = '28.614.48'
$ws.Range('E23').Value = '  +0.90%  '

$ws.Range('D24').Value = '11.05'
$ws.Range('E24').Value = '  -0.82%  '

$ws.Range('D25').Value = '2.234'
$ws.Range('E25').Value = '  -0.91%  '

$ws.Range('D26').Value = '2.088.82'
$ws.Range('E26').Value = '  +0.21%  '

$ws.Range('D27').Value = '161.56'
$ws.Range('E27').Value = '  +0.72%  '

$ws.Range('D28').Value = '20.64'
$ws.Range('E28').Value = '  -0.18%  '

$ws.Range('D29').Value = '2.330'
$ws.Range('E29').Value = '  -4.59%  '

$ws.Range('D30').Value = '125.79'
$ws.Range('E30').Value = '  -0.79%  '

$ws.Range('D31').Value = '0.1041'
$ws.Range('E31').Value = '  -1.71%  '

$ws.Range('D32').Value = '1.038'
$ws.Range('E32').Value = '  -0.54%  '

$ws.Range('D33').Value = '5.776'
$ws.Range('E33').Value = '  -2.09%  '

$ws.Range('D34').Value = '3.615'
$ws.Range('E34').Value = '  -0.25%  '

$ws.Range('D35').Value = '0.02451'
$ws.Range('E35').Value = '  +0.42%  '

$ws.Range('D36').Value = '0.06515'
$ws.Range('E36').Value = '  -0.40%  '

$ws.Range('D37').Value = '0.2162'
$ws.Range('E37').Value = '  -1.18%  '

$ws.Range('D38').Value = '8.891'
$ws.Range('E38').Value = '  -3.99%  '

$ws.Range('D39').Value = '1.256'
$ws.Range('E39').Value = '  +0.77%  '

$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '5.043'
$ws.Range('E40').Value = '  +1.31%  '

$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '1.187'
$ws.Range('E41').Value = '  +0.25%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6411'
$ws.Range('E42').Value = '  -0.87%  '

$ws.Range('D43').Value = '11.06'
$ws.Range('E43').Value = '  -0.92%  '

$ws.Range('D44').Value = '1.010'
$ws.Range('E44').Value = '  +0.18%  '

$ws.Range('D45').Value = '0.6004'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('D46').Value = '13.02'
$ws.Range('E46').Value = '  -0.34%  '

$ws.Range('D47').Value = '3.691'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').Value = '2.000'
$ws.Range('E48').Value = '  -0.54%  '

$ws.Range('D49').Value = '1.213'
$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('D50').Value = '121.73'
$ws.Range('E50').Value = '  -0.12%  '

$ws.Range('D51').Value = '1.177'
$ws.Range('E51').Value = '  -7.95%  '
